$wb = $excel.ActiveWorkbook

$wsProbes = $wb.Worksheets.Item("Probes")
$wsEPS    = $wb.Worksheets.Item("EPS")

# ---------------------------------------------------------------------
# Probes sheet ("verified that probe config has no impact") - fill in
# a full probe-configuration table (p1..p11) that mirrors the TTC / EPS
# layout (name, face1, face2, offset1, offset2, mass, area, Ix, Iy, Iz,
# gimb arm).
# ---------------------------------------------------------------------

# Header row
$wsProbes.Range("A1").Value = "name"
$wsProbes.Range("B1").Value = "face1"
$wsProbes.Range("C1").Value = "face2"
$wsProbes.Range("D1").Value = "offset1"
$wsProbes.Range("E1").Value = "offset2"
$wsProbes.Range("F1").Value = "mass"
$wsProbes.Range("G1").Value = "area"
$wsProbes.Range("H1").Value = "Ix"
$wsProbes.Range("I1").Value = "Iy"
$wsProbes.Range("J1").Value = "Iz"
$wsProbes.Range("K1").Value = "gimb arm"

# Probe names p1..p11 (rows 2..12)
for ($i = 1; $i -le 11; $i++) {
    $wsProbes.Range("A$($i + 1)").Value = "p$i"
}

# Every probe sits on the same face for now
for ($r = 2; $r -le 12; $r++) {
    $wsProbes.Range("B$r").Value = "z+"
}

# Mass is the same placeholder (210) for every probe
for ($r = 2; $r -le 12; $r++) {
    $wsProbes.Range("F$r").Value = 210
}

# Offset chain down column D: pairs of probes share an offset, then it
# steps by 2.5 every second pair.
$wsProbes.Range("D2").Value = 1.2
$wsProbes.Range("D3").Formula = "=D2"
$wsProbes.Range("D4").Formula = "=D3"
$wsProbes.Range("D5").Formula = "=D4"
$wsProbes.Range("D6").Formula = "=D5+2.5"
$wsProbes.Range("D7").Formula = "=D6"
$wsProbes.Range("D8").Formula = "=D7"
$wsProbes.Range("D9").Formula = "=D8"
$wsProbes.Range("D10").Formula = "=D9+2.5"
$wsProbes.Range("D11").Formula = "=D10"
$wsProbes.Range("D12").Formula = "=D11"

# Row 2/3 carry the full set of derived columns (area / inertia chain),
# mirroring the EPS sheet's G/H/I/J/K formulas.
$wsProbes.Range("G2").Formula = "='Calculations Rough Input'!D3/2"
$wsProbes.Range("G3").Formula = "=G2"

$wsProbes.Range("H2").Value = 0
$wsProbes.Range("H3").Formula = "=H2"

$wsProbes.Range("I2").Formula = "=H2"
$wsProbes.Range("I3").Formula = "=I2"

$wsProbes.Range("J2").Formula = "=I2"
$wsProbes.Range("J3").Formula = "=J2"

$wsProbes.Range("K3").Formula = "=K2"

# Re-stamp (no-op) formatting over the used area so the otherwise-blank
# cells (face2/offset2/gimb-arm on rows 2-3, etc.) still materialize as
# real, styled cells instead of being left out of the sheet entirely.
$wsProbes.Range("A1:K3").Font.Bold = $false
$wsProbes.Range("B4:B12").Font.Bold = $false
$wsProbes.Range("D4:D5").Font.Bold = $false
$wsProbes.Range("D8:D9").Font.Bold = $false
$wsProbes.Range("D11:D12").Font.Bold = $false
$wsProbes.Range("F4:F12").Font.Bold = $false

# ---------------------------------------------------------------------
# Selections: EPS keeps a range selection from reviewing the sheet,
# Probes ends up the active tab with E17 selected/highlighted.
# ---------------------------------------------------------------------
$wsEPS.Activate()
$wsEPS.Range("A1:K3").Select()

$wsProbes.Activate()
$wsProbes.Range("E17").Select()
